# Auto-generated Excel COM-interop script
# Applies the 2024-10-18 data update to violent-crime-ytd workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5882
$ws.Range("K3").Value = 6057
$ws.Range("G4").Value = 1222
$ws.Range("K4").Value = 1260
$ws.Range("K5").Value = 431
$ws.Range("K6").Value = 6663
$ws.Range("G7").Value = 20026
$ws.Range("K7").Value = 20293

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 589
$ws.Range("K8").Value = 1343
$ws.Range("K9").Value = 90
$ws.Range("K11").Value = 385
$ws.Range("K12").Value = 36
$ws.Range("K14").Value = 104
$ws.Range("K19").Value = 584
$ws.Range("K20").Value = 478
$ws.Range("K27").Value = 189
$ws.Range("K29").Value = 1107
$ws.Range("K30").Value = 77
$ws.Range("K31").Value = 226
$ws.Range("K33").Value = 875
$ws.Range("K36").Value = 258
$ws.Range("K37").Value = 687
$ws.Range("K41").Value = 141
$ws.Range("K42").Value = 754
$ws.Range("K45").Value = 26
$ws.Range("K47").Value = 141
$ws.Range("K49").Value = 110
$ws.Range("K51").Value = 259
$ws.Range("K52").Value = 536
$ws.Range("K53").Value = 259
$ws.Range("K54").Value = 394
$ws.Range("K55").Value = 223
$ws.Range("G63").Value = 198
$ws.Range("K65").Value = 473
$ws.Range("K67").Value = 795
$ws.Range("K68").Value = 54
$ws.Range("K72").Value = 97
$ws.Range("K73").Value = 179
$ws.Range("K76").Value = 276
$ws.Range("K77").Value = 143
$ws.Range("K78").Value = 228
$ws.Range("K79").Value = 502
$ws.Range("K80").Value = 71
$ws.Range("K81").Value = 15
$ws.Range("K83").Value = 453
$ws.Range("K84").Value = 160
$ws.Range("K85").Value = 949
$ws.Range("K88").Value = 219
$ws.Range("K90").Value = 187
$ws.Range("K91").Value = 230
$ws.Range("K94").Value = 272
$ws.Range("K97").Value = 161
$ws.Range("K98").Value = 97
$ws.Range("G101").Value = 20026
$ws.Range("K101").Value = 20293

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 197
$ws.Range("K3").Value = 192
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 589

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 137
$ws.Range("K3").Value = 99
$ws.Range("K7").Value = 385

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 309
$ws.Range("K3").Value = 325
$ws.Range("K6").Value = 235
$ws.Range("K7").Value = 949

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 144
$ws.Range("K3").Value = 156
$ws.Range("K5").Value = 19
$ws.Range("K7").Value = 536

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 66
$ws.Range("K3").Value = 68
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 372
$ws.Range("K3").Value = 411
$ws.Range("K5").Value = 39
$ws.Range("K7").Value = 1343

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 161
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 453

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K4").Value = 42
$ws.Range("K5").Value = 22
$ws.Range("K7").Value = 875

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 226
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 687

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 153
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 175
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 222
$ws.Range("K3").Value = 283
$ws.Range("K6").Value = 227
$ws.Range("K7").Value = 795

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 212
$ws.Range("K7").Value = 394

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 399
$ws.Range("K6").Value = 313
$ws.Range("K7").Value = 1107

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 179
$ws.Range("K7").Value = 584

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 205
$ws.Range("K3").Value = 231
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 280
$ws.Range("K7").Value = 754

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 223

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 170
$ws.Range("K3").Value = 162
$ws.Range("K7").Value = 502

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 159
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 478

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 102
$ws.Range("K7").Value = 258

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 72
$ws.Range("K6").Value = 121
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 55
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 68
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 25
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 15
